# "incluí alguns caracteres de teste" - set a test value in A1 and leave
# the selection on Q9, matching the sheet1.xml produced by the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the test value into A1 (becomes <c r="A1"><v>2</v></c>).
$ws.Range("A1").Value = 2

# Move/leave the active selection on Q9 (becomes <selection activeCell="Q9" sqref="Q9"/>).
$ws.Range("Q9").Select()
